$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 56 (pushes old rows 56..165 down to 57..166,
# and extends the sheet's used range / dimension to R166).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with its data. The "constant" columns
# (A,B,C,E,F,G,H,I,N,O,Q,R) are identical for every data row in this sheet,
# so copy them straight from the row that used to be 56 (now shifted to 57).
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44469
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112037
$ws.Range("G56").Value = "Cebollín"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 90
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = 6000
$ws.Range("N56").Value = "$/paquete 36 unidades"
$ws.Range("O56").Value = "Región Metropolitana"
$ws.Range("P56").Value = 167
$ws.Range("Q56").Value = 36
$ws.Range("R56").Value = "Hortaliza"
